# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the existing "Late" / "Outstanding" columns one to the
# right, and move the active selection to S8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (14th column) - this shifts
# N:P -> O:Q and copies formatting from the column to the left (M).
$ws.Columns("N").Insert()

# The inserted column picks up the width of its left neighbour (column M,
# "In Advance"), matching the source workbook.
$ws.Columns("N").ColumnWidth = 10.7109375

# Match the post-edit UI selection recorded in the workbook.
$ws.Range("S8").Select()
